# Export with no is_pref and no lev distance
# Rewrites the id/speaker_variant rows into their new export order and
# blanks out the "is_prefered" (column D) flag for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2-7: column B (id), column C (speaker_variant).
# Column D ("is_prefered") no longer gets an "x" marker for any row.
$rows = @(
    @{ Row = 2; Id = "#fem";        Variant = "Fem" },
    @{ Row = 3; Id = "#ouwekennis"; Variant = "Ouwekennis" },
    @{ Row = 4; Id = "#pefroen";    Variant = "Pefroen" },
    @{ Row = 5; Id = "#lysje";      Variant = "Lysje" },
    @{ Row = 6; Id = "#otje";       Variant = "Otje" },
    @{ Row = 7; Id = "#ritzaart";   Variant = "Ritzaart" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Id
    $ws.Cells.Item($r.Row, 3).Value = $r.Variant
    $ws.Cells.Item($r.Row, 4).Value = ""
}
